$d = $word.ActiveDocument

# Insert two new paragraphs at the very start of the document body:
#   1) "We have addressed the editor's comment (about spiking neurons) in the discussion."
#   2) an empty paragraph (Times New Roman formatting carried on the paragraph mark)
# Both use the "PreformattedText" style, matching the surrounding document.

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PreformattedText"/>
              <w:rPr/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
              </w:rPr>
              <w:t>We have addressed the editor’s comment (about spiking neurons) in the discussion.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="PreformattedText"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr/>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r = $d.Range(0, 0)
$r.InsertXML($xml)
